# Revision: add a new "chemical_recycling_pyrolysis" parameter row right
# after the existing "chemical_recycling_gasification" row, pushing every
# row below it down by one (dimension grows from A1:C24 to A1:C25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "chemical_recycling_gasification" lives on row 9; insert a fresh row at
# row 10 (all rows 10-24 shift down to 11-25) and populate it.
$ws.Rows(10).Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
